# Weekly crypto price/volume refresh (GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage for cells whose
# content would otherwise be auto-detected as a number by Excel
# (e.g. "243.93"), matching the source data which is stored as text.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2
$ws.Range("D2").Value = "29.130.92"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").Value = "1.840.95"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
Set-TextValue "D4" "0.9995"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
Set-TextValue "D5" "243.93"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
Set-TextValue "D6" "0.6262"
$ws.Range("E6").Value = "  -1.19%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
Set-TextValue "D8" "0.07551"
$ws.Range("E8").Value = "  -0.48%  "

# Row 9
Set-TextValue "D9" "0.2949"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
Set-TextValue "D10" "23.38"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11
Set-TextValue "D11" "0.07709"
$ws.Range("E11").Value = "  -0.41%  "

# Row 12
$ws.Range("D12").Value = "1.842.03"
$ws.Range("E12").Value = "  +0.44%  "

# Row 13
Set-TextValue "D13" "5.031"
$ws.Range("E13").Value = "  +0.80%  "

# Row 14
Set-TextValue "D14" "0.6783"
$ws.Range("E14").Value = "  +1.10%  "

# Row 15
Set-TextValue "D15" "83.22"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
Set-TextValue "D16" "0.000009331"
$ws.Range("E16").Value = "  -5.25%  "

# Row 17
Set-TextValue "D17" "5.988"
$ws.Range("E17").Value = "  -2.05%  "

# Row 18
$ws.Range("D18").Value = "29.125.90"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").Value = "2.080.78"
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
$ws.Range("E20").Value = "  +1.26%  "

# Row 21
Set-TextValue "D21" "229.53"
$ws.Range("E21").Value = "  +1.17%  "

# Row 22
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
Set-TextValue "D23" "7.181"
$ws.Range("E23").Value = "  -0.50%  "

# Row 24
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
Set-TextValue "D25" "160.57"
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
Set-TextValue "D26" "0.1405"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
Set-TextValue "D27" "8.564"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
Set-TextValue "D29" "1.497"
$ws.Range("E29").Value = "  -0.45%  "

# Row 30
Set-TextValue "D30" "4.194"
$ws.Range("E30").Value = "  +1.82%  "

# Row 31
$ws.Range("E31").Value = "  +2.41%  "

# Row 32
Set-TextValue "D32" "0.05578"
$ws.Range("E32").Value = "  +3.20%  "

# Row 33
$ws.Range("E33").Value = "  +0.36%  "

# Row 34
Set-TextValue "D34" "0.7498"
$ws.Range("E34").Value = "  +0.42%  "

# Row 35
$ws.Range("E35").Value = "  -0.48%  "

# Row 36
$ws.Range("E36").Value = "  +0.66%  "

# Row 37
Set-TextValue "D37" "2.669"
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("D38").Value = "1.238.09"
$ws.Range("E38").Value = "  -0.39%  "

# Row 39
Set-TextValue "D39" "2.772"
$ws.Range("E39").Value = "  +0.48%  "

# Row 40
Set-TextValue "D40" "0.01792"
$ws.Range("E40").Value = "  -0.21%  "

# Row 41
Set-TextValue "D41" "6.611"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
Set-TextValue "D42" "0.8997"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
Set-TextValue "D44" "102.44"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "66.83"
$ws.Range("E45").Value = "  +3.18%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.984.52"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
Set-TextValue "D47" "0.00000000123"
$ws.Range("E47").Value = "  +0.04%  "

# Row 48
$ws.Range("E48").Value = "  -0.41%  "

# Row 49
Set-TextValue "D49" "0.4090"
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
Set-TextValue "D50" "9.102"
$ws.Range("E50").Value = "  +1.01%  "

# Row 51
Set-TextValue "D51" "0.07291"
$ws.Range("E51").Value = "  +17.48%  "
